# Generate Report for Handback
# The handback transform failed for the second source file in both the
# zh-cn and de-de worksheets: update the Status to reflect the failure and
# populate the Error Detail column with the diagnostic message, widening
# that column so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status column - row 3 is the "83905d60-..." file that failed handback.
# Update every cell that shows this status (Overview summary columns E/F as
# well as each language sheet's Status column) so the shared text changes
# everywhere it is displayed.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Error Detail column (P) - row 3.
$zhcn.Range("P3").Value = "Handback file name: daywybmf.l3h is different with handoff file name: 83905d60-f64b-4cf7-8654-c03c0af893c6.972a74a8b66380e9ec38a572b6a23901823eeada.zh-cn."
$dede.Range("P3").Value = "Handback file name: daywybmf.l3h is different with handoff file name: 83905d60-f64b-4cf7-8654-c03c0af893c6.972a74a8b66380e9ec38a572b6a23901823eeada.de-de."

# Widen the Error Detail column (P) to fit the new message (matches the
# width already used by the other wide columns on these sheets).
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
